$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.643.07"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.205.37"
$ws.Range("E3").Value = "  -2.30%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.62"
$ws.Range("E5").Value = "  -1.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  -4.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.29"
$ws.Range("E7").Value = "  -6.42%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.402"
$ws.Range("E9").Value = "  -2.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.12"
$ws.Range("E10").Value = "  -5.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("E11").Value = "  -2.44%  "

$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.532.16"
$ws.Range("E13").Value = "  -2.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.37"
$ws.Range("E14").Value = "  -5.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.22"
$ws.Range("E15").Value = "  -2.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.65"
$ws.Range("E16").Value = "  -0.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.793"
$ws.Range("E17").Value = "  -4.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.206.57"
$ws.Range("E18").Value = "  -2.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.612.20"
$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.03"
$ws.Range("E20").Value = "  -2.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0900"
$ws.Range("E21").Value = "  -4.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("E22").Value = "  -2.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.46"
$ws.Range("E23").Value = "  -4.37%  "

$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  -2.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").Value = "  -3.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.62"
$ws.Range("E27").Value = "  -2.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "169.19"
$ws.Range("E28").Value = "  -1.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.140"
$ws.Range("E29").Value = "  -6.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("E30").Value = "  +0.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.74"
$ws.Range("E31").Value = "  -3.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.58"
$ws.Range("E32").Value = "  -8.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.121"
$ws.Range("E33").Value = "  -3.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.01"
$ws.Range("E34").Value = "  -1.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.62"
$ws.Range("E35").Value = "  -3.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0647"
$ws.Range("E36").Value = "  +0.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.36"
$ws.Range("E37").Value = "  -5.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.33"
$ws.Range("E38").Value = "  -8.90%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.53"
$ws.Range("E39").Value = "  -8.40%  "

$ws.Range("E40").Value = "  -7.69%  "

$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0237"
$ws.Range("E42").Value = "  -2.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.54"
$ws.Range("E43").Value = "  -3.06%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.25"
$ws.Range("E44").Value = "  -5.56%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0950"
$ws.Range("E45").Value = "  -3.88%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.19"
$ws.Range("E46").Value = "  -3.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.38"
$ws.Range("E47").Value = "  -13.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.465.63"
$ws.Range("E48").Value = "  -3.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.34"
$ws.Range("E49").Value = "  -7.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.77"
$ws.Range("E50").Value = "  -1.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.07"
$ws.Range("E51").Value = "  -6.28%  "
